# Quarterly indexing esoteric bug-fix operation
#
# Each data row (rows 2-16) held a staircase-shaped series of QoQ errors
# starting in column B. The first value in the series was being skipped
# (an off-by-one on the quarterly index), so every row needs a freshly
# computed leading value inserted at column B, with the previously
# stored values shifting one column to the right (B->C, C->D, ... );
# the last stored value in each row's series (which has no "next"
# quarter to host it) is dropped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The newly computed (previously missing) leading QoQ error for each row.
$newLeadValues = @{
    2  = -0.2177157015159319
    3  = -0.1395947820665385
    4  = -0.3119065001142551
    5  = 0.7021231295320197
    6  = 1.514070997382048
    7  = 0.2163102553365951
    8  = 0.3684555432821496
    9  = 0.661541622456546
    10 = -0.07992401592518952
    11 = 0.1551026493581833
    12 = -0.08373363042288225
    13 = 0.1925427069667326
    14 = -0.4379379024501944
    15 = 0.2324016585002178
    16 = -0.09587373626955231
}

$firstCol = 2   # column B
$lastCol  = 11  # column K

for ($row = 2; $row -le 16; $row++) {

    # Determine how many values are currently populated in this row
    # (columns B.. up to the first empty cell / lastCol).
    $count = 0
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($row, $c)
        $val = $cell.Value()
        if ($val -eq $null) { break }
        $count = $count + 1
    }

    # Read the existing values (column B .. B+count-1) before overwriting
    # anything, then shift them one column to the right, dropping the
    # value that would fall past column K.
    $existing = @()
    for ($i = 0; $i -lt $count; $i++) {
        $existing += , $ws.Cells.Item($row, $firstCol + $i).Value()
    }

    $shiftCount = $count
    if (($firstCol + $shiftCount) -gt $lastCol) {
        $shiftCount = $lastCol - $firstCol
    }

    for ($i = $shiftCount - 1; $i -ge 0; $i--) {
        $destCol = $firstCol + $i + 1
        $ws.Cells.Item($row, $destCol).Value = $existing[$i]
    }

    # Write the newly derived leading value into column B.
    $ws.Cells.Item($row, $firstCol).Value = $newLeadValues[$row]
}
